$wb = $excel.ActiveWorkbook

# --- Add "CheckedIn" sheet (after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$checkedIn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$checkedIn.Name = "CheckedIn"
$checkedIn.Range("A1").Value = "visitor_NIC"
$checkedIn.Range("A1").Interior.Color = 65535
$checkedIn.Range("A2").Value = "6348445764v"
$checkedIn.Range("A3").Value = "7821459632v"
$checkedIn.Range("A4").Value = "4578963245v"
$checkedIn.Columns.Item(1).ColumnWidth = 13

# --- Add "Overdue" sheet (after CheckedIn) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$overdue = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$overdue.Name = "Overdue"
$overdue.Range("A1").Value = "visitor_NIC"
$overdue.Range("A1").Interior.Color = 65535
$overdue.Range("A2").Value = "5612345783v"
$overdue.Range("A3").Value = "8695748612v"
$overdue.Columns.Item(1).ColumnWidth = 12

# --- Add "manageBuilding" sheet (after Overdue) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$manageBuilding = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$manageBuilding.Name = "manageBuilding"
$manageBuilding.Range("A1").Value = "BuildingName"
$manageBuilding.Range("B1").Value = "floorName"
$manageBuilding.Range("A1:B1").Interior.Color = 65535
$manageBuilding.Range("A2").Value = "buildingNo1"
$manageBuilding.Range("A3").Value = "buildingNo2"
$manageBuilding.Range("B2").Value = "floorNo1"
$manageBuilding.Range("B3").Value = "floorNo2"
$manageBuilding.Columns.Item(1).ColumnWidth = 18.333333333333332
$manageBuilding.Columns.Item(2).ColumnWidth = 14.166666666666666

# manageBuilding ends up the active/selected tab, matching the target workbook view.
$manageBuilding.Activate()

Write-Output "done"
